# Update countries & provincias Spain
# - Tunez overtakes Bulgaria in the country ranking (rows 84/85 swap identity)
# - Refreshed case numbers for Estados Unidos (row 4), Alemania (row 8),
#   and the Tunez/Bulgaria rows (84/85)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 670598
$ws.Range("C4").Value = 22450
$ws.Range("D4").Value = 57232
$ws.Range("E4").Value = 579133
$ws.Range("G4").Value = 1645
$ws.Range("H4").Value = 34233

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 136569
$ws.Range("C8").Value = 1816
$ws.Range("E8").Value = 55626
$ws.Range("G8").Value = 139
$ws.Range("H8").Value = 3943

# --- Rows 84/85: Tunez now ranks above Bulgaria, so they swap rows ---
# Row 84 becomes Tunez with updated totals
$ws.Range("A84").Value = "Tunez"
$ws.Range("B84").Value = 822
$ws.Range("C84").Value = 42
$ws.Range("D84").Value = 43
$ws.Range("E84").Value = 742
$ws.Range("F84").Value = 89
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 37

# Row 85 becomes Bulgaria, keeping its previous totals
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 800
$ws.Range("C85").Value = 53
$ws.Range("D85").Value = 122
$ws.Range("E85").Value = 640
$ws.Range("F85").Value = 37
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 38
